# Append two new daily rows (2025-10-30 and 2025-10-31, 3-day lead time)
# to the flood trigger analysis sheet for the Tuguegarao Buntun Bridge
# station on the Cagayan River Basin (Philippines / G4630 sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style to reuse for the new "forecast_date" text cells so they end up
# with the same (default) formatting as the rest of that column instead
# of picking up whatever format was used to force them to stay text.
$defaultStyle = $ws.Range("H30").Style

$newRows = @(
    @("Philippines","philippines","Cagayan River Basin","cagayan","Tuguegarao Buntun Bridge","G4630","primary","2025-10-30",3,17.62499999999999,121.6749999999997,5,9742.726601479721,"LOW",6867.840370863779,9742.726601479721,50,0,0,1302.41015625,1322.48583984375,1131.9140625,1637.3828125,1221.78515625,1375.154296875,$false,-86.63197470765328),
    @("Philippines","philippines","Cagayan River Basin","cagayan","Tuguegarao Buntun Bridge","G4630","primary","2025-10-31",3,17.62499999999999,121.6749999999997,5,9742.726601479721,"LOW",6867.840370863779,9742.726601479721,50,0,0,1289.8828125,1303.508911132812,1129.0390625,1745.171875,1226.97265625,1342.3671875,$false,-86.76055620503512)
)

$startRow = 31

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i

    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $value = $rowValues[$col - 1]

        if ($col -eq 8) {
            # forecast_date (column H) looks like a date string
            # ("2025-10-30"); force it to stay plain text instead of
            # letting Excel auto-convert it to a date serial, then
            # restore the normal/default cell style.
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.Style = $defaultStyle
        } else {
            $cell.Value = $value
        }
    }
}
